# Apply cell updates per the commit diff (cryptos price/volume refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.908.91"
$ws.Range("E2").Value = "  +0.28%  "

$ws.Range("D3").Value = "2.084.47"
$ws.Range("E3").Value = "  +0.01%  "

$ws.Range("E4").Value = "  -0.01%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "233.41"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.42%  "

$ws.Range("E7").Value = "  +3.69%  "

$ws.Range("E8").Value = "  -0.01%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.395"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +2.21%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.0788"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +1.28%  "

$ws.Range("E11").Value = "  +1.47%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "14.77"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +2.81%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "21.23"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +0.91%  "

$ws.Range("E14").Value = "  +2.42%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "5.35"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +2.42%  "

$ws.Range("D16").Value = "2.086.21"
$ws.Range("E16").Value = "  +0.14%  "

$ws.Range("D17").Value = "37.805.46"
$ws.Range("E17").Value = "  +0.19%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "6.13"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -0.14%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "71.70"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +1.24%  "

$ws.Range("D20").Value = "0.0₃0849"
$ws.Range("E20").Value = "  +3.51%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "228.14"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +0.12%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.07%  "

$ws.Range("E23").Value = "  -0.28%  "

$ws.Range("E24").Value = "  +1.17%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "172.17"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +1.40%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "9.20"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +3.29%  "

$ws.Range("E27").Value = "  -1.19%  "

$ws.Range("E28").Value = "  -1.42%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "19.50"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +0.32%  "

$ws.Range("E30").Value = "  +1.86%  "

$ws.Range("E31").Value = "  +2.81%  "

$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "4.72"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +3.22%  "

$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.0632"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +1.21%  "

$ws.Range("E34").Value = "  +0.35%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "3.42"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +1.35%  "

$ws.Range("E36").Value = "  -0.77%  "

$ws.Range("E37").Value = "  -0.13%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "5.42"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +0.38%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.0983"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -1.14%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "99.07"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +2.15%  "

$ws.Range("E41").Value = "  +2.55%  "

$ws.Range("E42").Value = "  -1.01%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "16.88"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +8.24%  "

$ws.Range("D44").Value = "1.446.51"
$ws.Range("E44").Value = "  -0.41%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "1.15"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -0.26%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "4.18"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +3.14%  "

$ws.Range("E47").Value = "  +1.06%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "7.38"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +0.36%  "

$ws.Range("E49").Value = "  +0.18%  "

$ws.Range("D50").Value = "2.277.33"
$ws.Range("E50").Value = "  -0.06%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "46.88"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +1.18%  "

